# Add another alt name column ("altname3") for DRC and CAR.
#
# The source sheet has columns: A=iso3, B=altname, C=altname2,
# D=formername, E=formername2. We insert a new column D ("altname3")
# which pushes the existing formername/formername2 columns one slot to
# the right (D->E, E->F), then populate the new column for the two rows
# that gained an extra alternate name:
#   - COD (Democratic Republic of the Congo) -> "DRC"
#   - CAF (Central African Republic)          -> "CAR"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank column at D; existing D/E shift right to E/F.
$ws.Columns("D:D").Insert()

# Match the look of the neighbouring altname/altname2 columns.
$ws.Columns("D:D").ColumnWidth = 24.67

# New column header.
$ws.Range("D1").Value = "altname3"

# Populate the new alt names. Write DRC (row 49, COD) before CAR (row 34,
# CAF) so the shared-string table gets the same ordering as the source
# workbook (altname3, DRC, CAR).
$ws.Range("D49").Value = "DRC"
$ws.Range("B34").Value = "CAR"

# Leave the view focused on the newly-edited cell, like the original edit.
$ws.Range("B34").Select()
